# Update "想去人数" (number of people interested) values for several events
# on both the "展览" (Exhibition) and "全部类型" (All types) sheets, per
# updated bilibili show data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of event name (column C) -> new value for column F, applied to every
# worksheet that contains a matching row. Using the event name instead of a
# hard-coded row number keeps this correct even though the two sheets list
# the same events on different row numbers (the "全部类型" sheet has one
# extra row for a concert that isn't on the "展览" sheet).
$updates = @{
    "南宁·2024年中国(华南)动漫超级订货会暨动漫实业发展大会" = 40
    "南宁·原x穹x崩only" = 311
    "南宁·AP动漫游戏嘉年华" = 2719
    "南宁·布谷鸟动漫展4th" = 1910
    "南宁·恋与深空only" = 368
    "南宁·小蜜蜂动漫嘉年华2.0" = 118
    "南宁·AB动漫游戏嘉年华" = 962
}

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
